$wb = $excel.ActiveWorkbook

# --- Uitvallers sheet: append two dropped riders ---
$wsUit = $wb.Worksheets.Item("Uitvallers")
$wsUit.Range("A4").Value = "Jasper Philipsen"
$wsUit.Range("A5").Value = "Stefan Bissegger"
$wsUit.Range("A5").Select()

# --- Huidig sheet: fill in Etappe-4 (column E) winners for rows 6-19 ---
$wsHuidig = $wb.Worksheets.Item("Huidig")
$wsHuidig.Activate()

$wsHuidig.Range("E6").Value = "Tadej Pogacar"
$wsHuidig.Range("E7").Value = "Mathieu Van Der Poel"
$wsHuidig.Range("E8").Value = "Jonas Vingegaard"
$wsHuidig.Range("E9").Value = "Oscar Onley"
$wsHuidig.Range("E10").Value = "Romain Gregoire"
$wsHuidig.Range("E11").Value = "Joao Almeida"
$wsHuidig.Range("E12").Value = "Remco Evenepoel"
$wsHuidig.Range("E13").Value = "Matteo Jorgenson"
$wsHuidig.Range("E14").Value = "Mattias Skjelmose"
$wsHuidig.Range("E15").Value = "Kévin Vauquelin"
$wsHuidig.Range("E16").Value = "Mathieu Van Der Poel"
$wsHuidig.Range("E17").Value = "Jonathan Milan"
$wsHuidig.Range("E18").Value = "Tadej Pogacar"
$wsHuidig.Range("E19").Value = "Kévin Vauquelin"

$wsHuidig.Range("E19").Select()
